$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Questions Tracker")

# --- New "Strings" topic section, mirroring the existing section headers
# (ARRAYS / Binary Search / Sorting Algorithm in rows 16, 30 and 42). Copy
# the formatting of the row 42 header band down onto row 51 so it picks up
# the identical bold/fill/border style rather than minting new ones.
$ws.Range("A42:H42").Copy()
$ws.Range("A51:H51").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("A51:H51").Merge()
$ws.Range("A51").Value = "Strings"

# --- Row 52: "Length of String" entry
$ws.Range("B52").Value = "Length of String"
$ws.Range("D52").Value = 1
$ws.Range("E52").Value = "20/11/2022"
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = "No"
$ws.Range("H52").Value = "Basic Question"

# --- Row 53: "Reverse of string" entry (topic only, rest left blank for now)
$ws.Range("B53").Value = "Reverse of string"

# Reflect where the user's selection ended up after the edits
$ws.Range("C53").Select()
